# ---------------------------------------------------------------------------
# Edit described by the diff:
#   1. Add a new worksheet "ODI Batting Extra" (sheetId 4) at the end of the
#      workbook with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#      PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.
#   2. On the "ODI Batting" sheet, remove the stray empty B12/B13/B16 cells
#      and normalize the E12/E16 "blank" cells to a single space.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the "ODI Batting Extra" worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $cell = $extraSheet.Cells.Item(1, $col + 1)
    $cell.Value = $headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("3005", "2", "3", "0", "33.00%", "NO"),
    @("3015", "2", "2", "0", "5.00%", "NO"),
    @("3017", "2", "5", "0", "15.05%", "NO"),
    @("3018", "", "", "", "", "NO"),
    @("3019", "2", "8", "0", "29.94%", "NO"),
    @("3023", "", "", "", "", "NO"),
    @("3027", "", "", "", "", "NO"),
    @("3029", "2", "1", "0", "3.42%", "NO"),
    @("3032", "2", "5", "0", "14.01%", "NO"),
    @("4284", "5", "1", "0", "4.02%", "NO"),
    @("4287", "", "", "", "", "NO"),
    @("4294", "", "", "", "", "NO"),
    @("4297", "7", "0", "0", "4.99%", "NO"),
    @("4401", "", "", "", "", "NO"),
    @("4405", "", "", "", "", "NO"),
    @("4408", "", "", "", "", "NO")
)

$rowNum = 2
foreach ($r in $rows) {

    # A: MATCH_CODE - keep as literal text (e.g. "3005")
    $aCell = $extraSheet.Cells.Item($rowNum, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $r[0]
    $aCell.Style = "Normal"

    # B: BATTING_POSITION - a real number when present, blank otherwise
    if ($r[1] -ne "") {
        $extraSheet.Cells.Item($rowNum, 2).Value = [double]$r[1]
    }

    # C: NUM_4 - literal text
    if ($r[2] -ne "") {
        $cCell = $extraSheet.Cells.Item($rowNum, 3)
        $cCell.NumberFormat = "@"
        $cCell.Value = $r[2]
        $cCell.Style = "Normal"
    }

    # D: NUM_6 - literal text
    if ($r[3] -ne "") {
        $dCell = $extraSheet.Cells.Item($rowNum, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $r[3]
        $dCell.Style = "Normal"
    }

    # E: PERCENT_RUNS_OF_TOTAL - literal text (keep the trailing "%")
    if ($r[4] -ne "") {
        $eCell = $extraSheet.Cells.Item($rowNum, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $r[4]
        $eCell.Style = "Normal"
    }

    # F: MAN_OF_MATCH - plain text
    $extraSheet.Cells.Item($rowNum, 6).Value = $r[5]

    $rowNum = $rowNum + 1
}

# ---------------------------------------------------------------------------
# 2. Tidy up the "ODI Batting" sheet.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Remove the stray empty INNING_NUMBER cells.
$battingSheet.Range("B12").ClearContents()
$battingSheet.Range("B13").ClearContents()
$battingSheet.Range("B16").ClearContents()

# Normalize the blank MATCH_INNING cells to a single space.
$battingSheet.Range("E12").Value = " "
$battingSheet.Range("E16").Value = " "
